$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 357
$ws.Range("I2").Value = 977
$ws.Range("J2").Value = 3994
$ws.Range("K2").Value = 14
$ws.Range("L2").Value = 1097
$ws.Range("M2").Value = 61
$ws.Range("N2").Value = 743
$ws.Range("P2").Value = 20
$ws.Range("Q2").Value = 9
$ws.Range("R2").Value = 41
$ws.Range("S2").Value = 409
$ws.Range("T2").Value = 704
$ws.Range("V2").Value = 6324
$ws.Range("X2").Value = 6216
$ws.Range("Y2").Value = 7
$ws.Range("Z2").Value = 101
$ws.Range("AA2").Value = 46
